$wb = $excel.ActiveWorkbook

# Sheet 1: Overview
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511

# Sheet 2: zh-cn
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511

# Sheet 3: de-de
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
